$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Clear the stale "PASS" result first so its shared-string slot can be reclaimed ---
$ws.Range("E12").Value = ""

# --- Rows 5-8: Sanity Runmode (column D) flips from Y to N ---
$ws.Range("D5:D8").Value = "N"

# --- Add two new rows (13 and 14) below row 12, copying its look & feel ---
$ws.Range("A12:F12").Copy()
$ws.Range("A13:F14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 30

# Set the new TCID labels before touching Sprint, so shared strings append in
# the same order as the source edit (RequestForInformation, IssuedForApproval, Sprint3)
$ws.Range("A13").Value = "FLD_Transmittals_ActionRequired_CaC_RequestForInformation"
$ws.Range("B13").Value = "Verifies the Close and Cancel option for the user in the Transmital record"
$ws.Range("C13").Value = "N"
$ws.Range("D13").Value = "Y"
$ws.Range("E13").Value = ""

$ws.Range("A14").Value = "FLD_Transmittals_ActionRequired_CaC_IssuedForApproval"
$ws.Range("B14").Value = "Verifies the Close and Cancel option for the user in the Transmital record"
$ws.Range("C14").Value = "N"
$ws.Range("D14").Value = "Y"
$ws.Range("E14").Value = ""

# --- Row 12 Sprint moves from Sprint2 to Sprint3, and the new rows inherit Sprint3 too ---
$ws.Range("F12").Value = "Sprint3"
$ws.Range("F13").Value = "Sprint3"
$ws.Range("F14").Value = "Sprint3"

# --- Column D got a touch narrower after the edits ---
$ws.Columns.Item(4).ColumnWidth = 14.7

# --- Data validation ranges need to grow from row 12 to row 14 ---
$ws.Range("C2:D12").Validation.Delete()
$ws.Range("F2:F12").Validation.Delete()
$ws.Range("C2:D14").Validation.Add(3, 1, 1, """Y,N""")
$ws.Range("F2:F14").Validation.Add(3, 1, 1, """Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10""")

# --- Reset selection back to the top-left cell ---
$ws.Range("A1").Select() | Out-Null
